$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 116411.555
$ws.Range("I125").Value = 4210.6665
$ws.Range("K125").Value = 37895.9985
$ws.Range("M125").Value = -35435.9985
$ws.Range("H132").Value = 3291587.2
$ws.Range("I132").Value = 4050099.8
$ws.Range("J132").Value = 4700
$ws.Range("K132").Value = 12150299.4
$ws.Range("L132").Value = 14100
$ws.Range("M132").Value = -12147769.4
$ws.Range("N132").Value = -19160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1834.25
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""
$ws.Range("H36").Value = 4299.75
$ws.Range("I36").Value = 2399.6667
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 2399.6667
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -2053.6667
$ws.Range("N36").Value = -10692
$ws.Range("H39").Value = 1016
$ws.Range("I39").Value = 1016
$ws.Range("K39").Value = 1016
$ws.Range("M39").Value = -496
$ws.Range("H40").Value = 26007.75
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 26007.75
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 26007.75
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -26359.75
$ws.Range("H53").Value = 24543
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 24543
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 24543
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = -25907
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = ""
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H88").Value = 1840
$ws.Range("I88").Value = 1840
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1840
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1434
$ws.Range("N88").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H91").Value = 1840
$ws.Range("I91").Value = 1840
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1840
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -436
$ws.Range("N91").Value = ""
$ws.Range("H133").Value = 31983.334
$ws.Range("J133").Value = 31983.334
$ws.Range("L133").Value = 31983.334
$ws.Range("N133").Value = -37043.334
$ws.Range("H135").Value = 34758.43
$ws.Range("J135").Value = 34758.43
$ws.Range("L135").Value = 34758.43
$ws.Range("N135").Value = -44898.43

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 11925.546
$ws.Range("I36").Value = 888
$ws.Range("J36").Value = 21123.5
$ws.Range("K36").Value = 888
$ws.Range("L36").Value = 21123.5
$ws.Range("M36").Value = -354
$ws.Range("N36").Value = -22191.5
$ws.Range("H38").Value = 19500
$ws.Range("J38").Value = 19500
$ws.Range("L38").Value = 19500
$ws.Range("N38").Value = -20332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 5675.727
$ws.Range("I35").Value = 656.25
$ws.Range("J35").Value = 8544
$ws.Range("K35").Value = 656.25
$ws.Range("L35").Value = 8544
$ws.Range("M35").Value = -362.25
$ws.Range("N35").Value = -9132
$ws.Range("H38").Value = 21991.8
$ws.Range("I38").Value = 4639.3335
$ws.Range("J38").Value = 29428.572
$ws.Range("K38").Value = 4639.3335
$ws.Range("L38").Value = 29428.572
$ws.Range("M38").Value = -4262.3335
$ws.Range("N38").Value = -30182.572
$ws.Range("H46").Value = 21991.8
$ws.Range("I46").Value = 4639.3335
$ws.Range("J46").Value = 29428.572
$ws.Range("K46").Value = 4639.3335
$ws.Range("L46").Value = 29428.572
$ws.Range("M46").Value = -4428.3335
$ws.Range("N46").Value = -29850.572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3447.3333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3447.3333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 10341.9999
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -12213.9999
$ws.Range("H81").Value = 2402.3333
$ws.Range("I81").Value = 1013
$ws.Range("J81").Value = 2528.6365
$ws.Range("K81").Value = 3039
$ws.Range("L81").Value = 7585.9095
$ws.Range("M81").Value = -1916
$ws.Range("N81").Value = -9831.9095
$ws.Range("H82").Value = 5413.273
$ws.Range("I82").Value = 2019
$ws.Range("J82").Value = 7352.857
$ws.Range("K82").Value = 6057
$ws.Range("L82").Value = 22058.571
$ws.Range("M82").Value = -5651
$ws.Range("N82").Value = -22870.571
$ws.Range("H83").Value = 3447.3333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3447.3333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 31025.9997
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -40385.9997
$ws.Range("H84").Value = 2402.3333
$ws.Range("I84").Value = 1013
$ws.Range("J84").Value = 2528.6365
$ws.Range("K84").Value = 9117
$ws.Range("L84").Value = 22757.7285
$ws.Range("M84").Value = -3501
$ws.Range("N84").Value = -33989.7285
$ws.Range("H85").Value = 5413.273
$ws.Range("I85").Value = 2019
$ws.Range("J85").Value = 7352.857
$ws.Range("K85").Value = 6057
$ws.Range("L85").Value = 22058.571
$ws.Range("M85").Value = -4653
$ws.Range("N85").Value = -24866.571
$ws.Range("H106").Value = 6521.1113
$ws.Range("J106").Value = 6521.1113
$ws.Range("L106").Value = 19563.3339
$ws.Range("N106").Value = -21455.3339
$ws.Range("H131").Value = 887.8415
$ws.Range("I131").Value = 379.25
$ws.Range("J131").Value = 1051.9032
$ws.Range("K131").Value = 1137.75
$ws.Range("L131").Value = 3155.7096
$ws.Range("M131").Value = 3902.25
$ws.Range("N131").Value = -13235.7096
$ws.Range("H137").Value = 15763297
$ws.Range("I137").Value = 27784690
$ws.Range("J137").Value = 1337626.5
$ws.Range("K137").Value = 83354070
$ws.Range("L137").Value = 4012879.5
$ws.Range("M137").Value = -83348970
$ws.Range("N137").Value = -4023079.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3283.3333
$ws.Range("J80").Value = 4366.6665
$ws.Range("L80").Value = 4366.6665
$ws.Range("N80").Value = -6362.6665
$ws.Range("H83").Value = 3283.3333
$ws.Range("J83").Value = 4366.6665
$ws.Range("L83").Value = 21833.3325
$ws.Range("N83").Value = -31817.3325
$ws.Range("H86").Value = 50000
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 50000
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1490.7188
$ws.Range("I82").Value = 1212.625
$ws.Range("J82").Value = 1768.8125
$ws.Range("K82").Value = 1212.625
$ws.Range("L82").Value = 1768.8125
$ws.Range("M82").Value = -851.625
$ws.Range("N82").Value = -2490.8125
$ws.Range("H85").Value = 1490.7188
$ws.Range("I85").Value = 1212.625
$ws.Range("J85").Value = 1768.8125
$ws.Range("K85").Value = 1212.625
$ws.Range("L85").Value = 1768.8125
$ws.Range("M85").Value = 35.375
$ws.Range("N85").Value = -4264.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3045.3076
$ws.Range("I62").Value = 2473.625
$ws.Range("J62").Value = 3960
$ws.Range("K62").Value = 2473.625
$ws.Range("L62").Value = 3960
$ws.Range("M62").Value = -1849.625
$ws.Range("N62").Value = -5208
$ws.Range("H65").Value = 3045.3076
$ws.Range("I65").Value = 2473.625
$ws.Range("J65").Value = 3960
$ws.Range("K65").Value = 12368.125
$ws.Range("L65").Value = 19800
$ws.Range("M65").Value = -9248.125
$ws.Range("N65").Value = -26040
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("K88").Value = 10000
$ws.Range("M88").Value = -9594
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("K91").Value = 10000
$ws.Range("M91").Value = -8596
